$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 2 (OTROS)
$ws.Range("D2").Value = 6728.52
$ws.Range("E2").Value = -6728.52

# Row 3 (PORCELANATO)
$ws.Range("D3").Value = 14580.38
$ws.Range("E3").Value = -857.0399999999991
$ws.Range("F3").Value = 1.062451269151679

# Row 4 (TOTAL)
$ws.Range("D4").Value = 21308.9
$ws.Range("E4").Value = -7585.559999999999
$ws.Range("F4").Value = 1.552748820622385
